# "explorer klan van links e recht lopen"
# Fill in the logboek entries for "week 3" rows 13 and 14 (new activities
# for the explorer feature) and move the selection to just F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 3")

# Row 13: 10:10 - 12:45, id 7, "bezig met explorer"
$ws.Range("C13").Value = 0.4236111111111111
$ws.Range("D13").Value = 0.53125
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "bezig met explorer"

# Row 14: 13:30 - 14:00, id 8, "Explorer kan nu van links en recht lopen"
$ws.Range("C14").Value = 0.5625
$ws.Range("D14").Value = 0.58333333333333337
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = "Explorer kan nu van links en recht lopen"

# Selection collapses from F13:F14 to just F14
$null = $ws.Range("F14").Select()
